$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.485167622566223
$ws.Range("B1").Value = 2.707624197006226
$ws.Range("C1").Value = 6.765398502349854
$ws.Range("D1").Value = 1.751126646995544
$ws.Range("E1").Value = 0.8974946141242981
